$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.752.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.27%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.413.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.88%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.82%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.64'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  +0.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.385'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.997.01'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.85%  '

$ws.Range("E13").Value = '  -0.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.401.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.15%  '

$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.785.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.550'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.43%  '

$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.187'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.23%  '

$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.65%  '

$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("E31").Value = '  +1.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.19'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.95'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '168.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.447.01'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.89%  '

$ws.Range("E38").Value = '  +0.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '28.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0753'
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = '  +0.96%  '

$ws.Range("E42").Value = '  +1.80%  '

$ws.Range("E43").Value = '  +0.64%  '

$ws.Range("E44").Value = '  +4.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.501.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("E50").Value = '  -3.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.206'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.37%  '
